$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new value would otherwise be parsed as a number
# (these are plain decimal price strings that must remain literal text, matching
# the source workbook where the whole Price column is stored as text)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Write the updated Price (D) and Volume(1h) (E) values
$ws.Range("D2").Value = '67.579.26'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '3.771.77'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("D5").Value = '596.02'
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").Value = '168.58'
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("D7").Value = '3.773.23'
$ws.Range("E7").Value = '  -1.75%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("D10").Value = '0.164'
$ws.Range("E10").Value = '  -1.38%  '
$ws.Range("D11").Value = '6.46'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").Value = '0.453'
$ws.Range("E12").Value = '  -1.30%  '
$ws.Range("D13").Value = '0.0000274'
$ws.Range("E13").Value = '  +2.38%  '
$ws.Range("D14").Value = '36.41'
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("D15").Value = '4.402.33'
$ws.Range("E15").Value = '  -1.84%  '
$ws.Range("D16").Value = '3.761.10'
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("D17").Value = '18.57'
$ws.Range("D18").Value = '67.517.56'
$ws.Range("E18").Value = '  -1.39%  '
$ws.Range("D19").Value = '7.17'
$ws.Range("E19").Value = '  -3.16%  '
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("D21").Value = '10.57'
$ws.Range("E21").Value = '  -4.46%  '
$ws.Range("D22").Value = '467.08'
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("D23").Value = '0.717'
$ws.Range("E23").Value = '  -2.46%  '
$ws.Range("D24").Value = '83.52'
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("E25").Value = '  -9.94%  '
$ws.Range("D26").Value = '2.20'
$ws.Range("E26").Value = '  -1.39%  '
$ws.Range("D27").Value = '12.13'
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("D28").Value = '10.26'
$ws.Range("E28").Value = '  +2.06%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '2.91'
$ws.Range("E30").Value = '  -2.14%  '
$ws.Range("D31").Value = '3.916.84'
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("E32").Value = '  -0.93%  '
$ws.Range("E33").Value = '  -3.61%  '
$ws.Range("E34").Value = '  -3.84%  '
$ws.Range("D35").Value = '9.11'
$ws.Range("E35").Value = '  -2.98%  '
$ws.Range("D36").Value = '3.731.41'
$ws.Range("E36").Value = '  -2.01%  '
$ws.Range("D37").Value = '3.79'
$ws.Range("E37").Value = '  +2.34%  '
$ws.Range("E38").Value = '  -1.23%  '
$ws.Range("E39").Value = '  -1.70%  '
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("D41").Value = '5.78'
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = '8.68'
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("E46").Value = '  -2.39%  '
$ws.Range("D47").Value = '45.83'
$ws.Range("E47").Value = '  -2.79%  '
$ws.Range("D48").Value = '395.67'
$ws.Range("E48").Value = '  -4.54%  '
$ws.Range("D49").Value = '0.000270'
$ws.Range("E49").Value = '  -7.15%  '
$ws.Range("D50").Value = '139.28'
$ws.Range("E50").Value = '  -1.61%  '
